$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A26").Value = "mattia spagnolli"
$ws.Range("B26").Value = "MATTEO PILATI | Pinguini Trentini"
$ws.Range("C26").Value = "Leonardo Viola | Shark Attack"
$ws.Range("D26").Value = "Andrea Conzatti | FC Savignano"
$ws.Range("E26").Value = "FEDERICO NICOLODI | U.S. Guarna"
$ws.Range("F26").Value = "Simone Schonsberg | I Magnifici"
